$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21 (AdminDashboard): fill in the previously-blank Description (B21). ---
$ws.Range("B21").Value = "Done"

# --- Row 22 (ReplayEvent): fill in Description (B22) and correct Runmode (C22) N->Y typo to N. ---
$ws.Range("B22").Value = "Changes where the event replay or not"
$ws.Range("C22").Value = "N"

# --- Seed rows 23 & 24 with the same cell formatting (border/fill/font) as the
#     existing data rows, without disturbing any values, by copying formats only
#     from row 20 (a representative already-formatted data row). ---
$ws.Range("A20:C20").Copy()
$ws.Range("A23:C24").PasteSpecial(-4122)

# --- Row 23 (new): ReplayReport test case. ---
$ws.Range("A23").Value = "ReplayReport"
$ws.Range("B23").Value = "Changes where the report replay or not"
$ws.Range("C23").Value = "N"

# --- Row 24 (new): ImportSIMs test case (Description intentionally left blank). ---
$ws.Range("A24").Value = "ImportSIMs"
$ws.Range("C24").Value = "Y"

# --- Restore the author's final selection/cursor position. ---
$ws.Range("B20").Select()
